$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Title shape on the (single) slide: "Homogeneous coordinates & pinhole camera model (p = K[R t]Q)"
$shp = $s.Shapes.Item(5)
$tr = $shp.TextFrame.TextRange

# Split the single run into four runs, rewriting the wording/formula
# "... pinhole camera model (p = K[R t]Q)"  ->  "... pinhole camera model (q = K[R t]Q = PQ)"
$r1 = $tr.Characters(1, 41)
$r1.Text = "Homogeneous coordinates & pinhole camera "

$r2 = $tr.Characters(42, 9)
$r2.Text = "model (q "

$r3 = $tr.Characters(51, 7)
$r3.Text = "= K[R t"

$r4 = $tr.Characters(58, 8)
$r4.Text = "]Q = PQ)"
